$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F7: GDP Nowcast value update
$ws.Range("F7").Value = 5.3721

# Row 28 - Mich NTM Inflation Exp
$ws.Range("N28").Value = [DateTime]"2025-12-01"
$ws.Range("Q28").Value = 4.2
$ws.Range("R28").Value = 4.5
$ws.Range("S28").Value = 4.6
$ws.Range("T28").Value = 4.7
$ws.Range("U28").Value = 4.8

# Row 29 - 5yr, 5yr Forward
$ws.Range("N29").Value = [DateTime]"2026-01-23"
$ws.Range("Q29").Value = 2.18
$ws.Range("R29").Value = 2.2
$ws.Range("S29").Value = 2.26
$ws.Range("T29").Value = 2.26
$ws.Range("U29").Value = 2.27

# Row 30 - 10yr TIPS
$ws.Range("N30").Value = [DateTime]"2026-01-23"
$ws.Range("Q30").Value = 2.32
$ws.Range("R30").Value = 2.31
$ws.Range("S30").Value = 2.34
$ws.Range("T30").Value = 2.33
$ws.Range("U30").Value = 2.33

# Row 47 - FFR
$ws.Range("N47").Value = [DateTime]"2026-01-22"

# Row 48 - 2y UST
$ws.Range("N48").Value = [DateTime]"2026-01-22"
$ws.Range("Q48").Value = 3.61
$ws.Range("R48").Value = 3.6
$ws.Range("S48").Value = 3.6
$ws.Range("T48").Value = 3.59
$ws.Range("U48").Value = 3.56

# Row 49 - 5y UST
$ws.Range("N49").Value = [DateTime]"2026-01-22"
$ws.Range("Q49").Value = 3.85
$ws.Range("R49").Value = 3.83
$ws.Range("S49").Value = 3.86
$ws.Range("T49").Value = 3.82
$ws.Range("U49").Value = 3.77

# Row 50 - 10y UST
$ws.Range("N50").Value = [DateTime]"2026-01-22"
$ws.Range("Q50").Value = 4.26
$ws.Range("R50").Value = 4.26
$ws.Range("S50").Value = 4.3
$ws.Range("T50").Value = 4.24
$ws.Range("U50").Value = 4.17

# Row 52 - BAA
$ws.Range("N52").Value = [DateTime]"2026-01-22"
$ws.Range("Q52").Value = 5.85
$ws.Range("R52").Value = 5.88
$ws.Range("S52").Value = 5.95
$ws.Range("T52").Value = 5.87
$ws.Range("U52").Value = 5.82
